$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.079492092132568
$ws.Range("B1").Value = 3.12166428565979
$ws.Range("C1").Value = 2.850465774536133
$ws.Range("D1").Value = 3.624174833297729
$ws.Range("E1").Value = 5.07529878616333
